$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$vt = [char]11

$c0 = $t.Cell(1, 1)
$c0.Range.Text = "18 x 79" + $vt + "  7    9" + $vt + "  ----" + $vt + "1|    |" + $vt + "8|    |"

$c1 = $t.Cell(1, 2)
$c1.Range.Text = "92 x 20" + $vt + "  2    0" + $vt + "  ----" + $vt + "9|    |" + $vt + "2|    |"

$c2 = $t.Cell(1, 3)
$c2.Range.Text = "31 x 83" + $vt + "  8    3" + $vt + "  ----" + $vt + "3|    |" + $vt + "1|    |"

$c3 = $t.Cell(2, 1)
$c3.Range.Text = "58 x 82" + $vt + "  8    2" + $vt + "  ----" + $vt + "5|    |" + $vt + "8|    |"

$c4 = $t.Cell(2, 2)
$c4.Range.Text = "75 x 78" + $vt + "  7    8" + $vt + "  ----" + $vt + "7|    |" + $vt + "5|    |"

$c5 = $t.Cell(2, 3)
$c5.Range.Text = "59 x 74" + $vt + "  7    4" + $vt + "  ----" + $vt + "5|    |" + $vt + "9|    |"

$c6 = $t.Cell(3, 1)
$c6.Range.Text = "80 x 59" + $vt + "  5    9" + $vt + "  ----" + $vt + "8|    |" + $vt + "0|    |"

$c7 = $t.Cell(3, 2)
$c7.Range.Text = "22 x 51" + $vt + "  5    1" + $vt + "  ----" + $vt + "2|    |" + $vt + "2|    |"

$c8 = $t.Cell(3, 3)
$c8.Range.Text = "58 x 60" + $vt + "  6    0" + $vt + "  ----" + $vt + "5|    |" + $vt + "8|    |"

$c9 = $t.Cell(4, 1)
$c9.Range.Text = "82 x 73" + $vt + "  7    3" + $vt + "  ----" + $vt + "8|    |" + $vt + "2|    |"

$c10 = $t.Cell(4, 2)
$c10.Range.Text = "77 x 36" + $vt + "  3    6" + $vt + "  ----" + $vt + "7|    |" + $vt + "7|    |"

$c11 = $t.Cell(4, 3)
$c11.Range.Text = "67 x 85" + $vt + "  8    5" + $vt + "  ----" + $vt + "6|    |" + $vt + "7|    |"

$c12 = $t.Cell(5, 1)
$c12.Range.Text = "67 x 95" + $vt + "  9    5" + $vt + "  ----" + $vt + "6|    |" + $vt + "7|    |"

$c13 = $t.Cell(5, 2)
$c13.Range.Text = "24 x 29" + $vt + "  2    9" + $vt + "  ----" + $vt + "2|    |" + $vt + "4|    |"

$c14 = $t.Cell(5, 3)
$c14.Range.Text = "67 x 77" + $vt + "  7    7" + $vt + "  ----" + $vt + "6|    |" + $vt + "7|    |"
